# Add a "Save" column (H) to the s_vals sheet.
# H1 gets the header "Save" (matching the formatting of the other header
# cells in row 1), and H2:H48 get a 0/1 indicator value for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy formatting from the neighboring header cell (G1)
# so it gets the same bold/bordered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for H2:H48
$saveValues = @(0,0,0,0,1,0,0,0,0,1,0,0,1,0,0,0,0,0,0,1,1,0,0,0,1,0,0,0,0,0,1,1,0,0,1,0,1,1,0,0,0,1,0,1,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
